$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# Each entry: row number, new Display text (Column C), new Definition/abbreviation code (Column D)
$updates = @(
    @{ Row = 2; Display = '1ª Dose'; Code = 'D1' }
    @{ Row = 3; Display = '2ª Dose'; Code = 'D2' }
    @{ Row = 4; Display = '3ª Dose'; Code = 'D3' }
    @{ Row = 5; Display = '4ª Dose'; Code = 'D4' }
    @{ Row = 6; Display = '5ª Dose'; Code = 'D5' }
    @{ Row = 7; Display = '1º Reforço'; Code = 'R1' }
    @{ Row = 8; Display = '2º Reforço'; Code = 'R2' }
    @{ Row = 9; Display = 'Dose'; Code = 'D' }
    @{ Row = 10; Display = 'Única'; Code = 'DU' }
    @{ Row = 11; Display = 'Revacinação'; Code = 'REV' }
    @{ Row = 12; Display = 'Tratamento com uma dose'; Code = 'T1' }
    @{ Row = 13; Display = 'Tratamento com duas doses'; Code = 'T2' }
    @{ Row = 14; Display = 'Tratamento com três doses'; Code = 'T3' }
    @{ Row = 15; Display = 'Tratamento com quatro doses'; Code = 'T4' }
    @{ Row = 16; Display = 'Tratamento com cinco doses'; Code = 'T5' }
    @{ Row = 17; Display = 'Tratamento com seis doses'; Code = 'T6' }
    @{ Row = 18; Display = 'Tratamento com sete doses'; Code = 'T7' }
    @{ Row = 19; Display = 'Tratamento com oito doses'; Code = 'T8' }
    @{ Row = 20; Display = 'Tratamento com nove doses'; Code = 'T9' }
    @{ Row = 21; Display = 'Tratamento com dez doses'; Code = 'T10' }
    @{ Row = 22; Display = 'Tratamento com onze doses'; Code = 'T11' }
    @{ Row = 23; Display = 'Tratamento com doze doses'; Code = 'T12' }
    @{ Row = 24; Display = 'Tratamento com treze doses'; Code = 'T13' }
    @{ Row = 25; Display = 'Tratamento com quartorze doses'; Code = 'T14' }
    @{ Row = 26; Display = 'Tratamento com quinze doses'; Code = 'T15' }
    @{ Row = 27; Display = 'Tratamento com dezesseis doses'; Code = 'T16' }
    @{ Row = 28; Display = 'Tratamento com dezessete doses'; Code = 'T17' }
    @{ Row = 29; Display = 'Tratamento com dezoito doses'; Code = 'T18' }
    @{ Row = 30; Display = 'Tratamento com dezenove doses'; Code = 'T19' }
    @{ Row = 31; Display = 'Tratamento com vinte doses'; Code = 'T20' }
    @{ Row = 32; Display = 'Tratamento com vinte e quatro doses'; Code = 'T24' }
    @{ Row = 33; Display = '1ª Dose Revacinação'; Code = 'D1REV' }
    @{ Row = 34; Display = '2ª Dose Revacinação'; Code = 'D2REV' }
    @{ Row = 35; Display = '3ª Dose Revacinação'; Code = 'D3REV' }
    @{ Row = 36; Display = '4ª Dose Revacinação'; Code = 'D4REV' }
    @{ Row = 37; Display = 'Dose Inicial'; Code = 'DI' }
    @{ Row = 38; Display = 'Dose Adicional'; Code = 'DA' }
    @{ Row = 39; Display = 'Reforço'; Code = 'REF' }
    @{ Row = 40; Display = '3º Reforço'; Code = 'R3' }
    @{ Row = 41; Display = '4º Reforço'; Code = 'R4' }
    @{ Row = 42; Display = '5º Reforço'; Code = 'R5' }
    @{ Row = 43; Display = '6º Reforço'; Code = 'R6' }
    @{ Row = 44; Display = '5ª Dose Revacinação'; Code = 'D5REV' }
    @{ Row = 45; Display = '1ª Dose Fracionada'; Code = 'D1F' }
    @{ Row = 46; Display = '2ª Dose Fracionada'; Code = 'D2F' }
    @{ Row = 47; Display = '3ª Dose Fracionada'; Code = 'D3F' }
    @{ Row = 48; Display = '4ª Dose Fracionada'; Code = 'D4F' }
    @{ Row = 49; Display = '5ª Dose Fracionada'; Code = 'D5F' }
    @{ Row = 50; Display = '1ª Dose Dobrada'; Code = 'D1D' }
    @{ Row = 51; Display = '2ª Dose Dobrada'; Code = 'D2D' }
    @{ Row = 52; Display = '3ª Dose Dobrada'; Code = 'D3D' }
    @{ Row = 53; Display = '4ª Dose Dobrada'; Code = 'D4D' }
    @{ Row = 54; Display = '1ª Dose Revacinação Dobrada'; Code = 'D1REVD' }
    @{ Row = 55; Display = '2ª Dose Revacinação Dobrada'; Code = 'D2REVD' }
    @{ Row = 56; Display = '3ª Dose Revacinação Dobrada'; Code = 'D3REVD' }
    @{ Row = 57; Display = '4ª Dose Revacinação Dobrada'; Code = 'D4REVD' }
    @{ Row = 58; Display = 'Dose Zero'; Code = 'D0' }
    @{ Row = 59; Display = 'Reforço Zero'; Code = 'R0' }
    @{ Row = 60; Display = 'Profilaxia/Tratamento com 1 frasco-ampola/ampola'; Code = 'P/T1' }
    @{ Row = 61; Display = 'Profilaxia/Tratamento com 2 frascos-ampolas/ampolas'; Code = 'P/T2' }
    @{ Row = 62; Display = 'Profilaxia/Tratamento com 3 frascos-ampolas/ampolas'; Code = 'P/T3' }
    @{ Row = 63; Display = 'Profilaxia/Tratamento com 4 frascos-ampolas/ampolas'; Code = 'P/T4' }
    @{ Row = 64; Display = 'Profilaxia/Tratamento com 5 frascos-ampolas/ampolas'; Code = 'P/T5' }
    @{ Row = 65; Display = 'Profilaxia/Tratamento com 6 frascos-ampolas/ampolas'; Code = 'P/T6' }
    @{ Row = 66; Display = 'Profilaxia/Tratamento com 7 frascos-ampolas/ampolas'; Code = 'P/T7' }
    @{ Row = 67; Display = 'Profilaxia/Tratamento com 8 frascos-ampolas/ampolas'; Code = 'P/T8' }
    @{ Row = 68; Display = 'Profilaxia/Tratamento com 9 frascos-ampolas/ampolas'; Code = 'P/T9' }
    @{ Row = 69; Display = 'Profilaxia/Tratamento com 10 frascos-ampolas/ampolas'; Code = 'P/T10' }
    @{ Row = 70; Display = 'Profilaxia/Tratamento com 11 frascos-ampolas/ampolas'; Code = 'P/T11' }
    @{ Row = 71; Display = 'Profilaxia/Tratamento com 12 frascos-ampolas/ampolas'; Code = 'P/T12' }
    @{ Row = 72; Display = 'Profilaxia/Tratamento com 13 frascos-ampolas/ampolas'; Code = 'P/T13' }
    @{ Row = 73; Display = 'Profilaxia/Tratamento com 14 frascos-ampolas/ampolas'; Code = 'P/T14' }
    @{ Row = 74; Display = 'Profilaxia/Tratamento com 15 frascos-ampolas/ampolas'; Code = 'P/T15' }
    @{ Row = 75; Display = 'Profilaxia/Tratamento com 16 frascos-ampolas/ampolas'; Code = 'P/T16' }
    @{ Row = 76; Display = 'Profilaxia/Tratamento com 17 frascos-ampolas/ampolas'; Code = 'P/T17' }
    @{ Row = 77; Display = 'Profilaxia/Tratamento com 18 frascos-ampolas/ampolas'; Code = 'P/T18' }
    @{ Row = 78; Display = 'Profilaxia/Tratamento com 19 frascos-ampolas/ampolas'; Code = 'P/T19' }
    @{ Row = 79; Display = 'Profilaxia/Tratamento com 20 frascos-ampolas/ampolas'; Code = 'P/T20' }
    @{ Row = 80; Display = 'Profilaxia/Tratamento com 21 frascos-ampolas/ampolas'; Code = 'P/T21' }
    @{ Row = 81; Display = 'Profilaxia/Tratamento com 22 frascos-ampolas/ampolas'; Code = 'P/T22' }
    @{ Row = 82; Display = 'Profilaxia/Tratamento com 23 frascos-ampolas/ampolas'; Code = 'P/T23' }
    @{ Row = 83; Display = 'Profilaxia/Tratamento com 24 frascos-ampolas/ampolas'; Code = 'P/T24' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value2 = $u.Display
    $ws.Cells.Item($u.Row, 4).Value2 = $u.Code
}
